# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# "Bad Drivers" table (row 3-4)
$ws.Range("C3").Value = 253
$ws.Range("D3").Value = 98.2
$ws.Range("C4").Value = 253

# "Good Drivers" table (rows 12-13)
$ws.Range("B12").Value = 11140
$ws.Range("B13").Value = 14487
